# Updates the cryptos worksheet with the latest scraped values.
# Only cells that actually changed (per the upstream diff) are touched;
# everything else is left exactly as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; D="41.700.32"; E="  +0.18%  "}
    @{Row=3; D="2.464.89"; E="  -0.10%  "}
    @{Row=4; E="  +0.14%  "}
    @{Row=5; D="319.42"; E="  +0.58%  "}
    @{Row=6; D="91.38"; E="  -1.00%  "}
    @{Row=7; D="0.548"; E="  -0.73%  "}
    @{Row=8; E="  +0.10%  "}
    @{Row=9; D="0.505"; E="  -2.00%  "}
    @{Row=10; B="Avalanche"; C="https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D="32.61"; E="  +0.06%  "}
    @{Row=11; B="Dogecoin"; C="https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; D="0.0847"; E="  -5.13%  "}
    @{Row=12; E="  -0.79%  "}
    @{Row=13; D="2.844.52"; E="  -0.12%  "}
    @{Row=14; D="6.83"; E="  -0.88%  "}
    @{Row=15; D="15.50"; E="  -0.64%  "}
    @{Row=16; D="2.452.04"; E="  -0.49%  "}
    @{Row=17; D="0.789"; E="  +0.76%  "}
    @{Row=18; D="41.583.24"; E="  -0.04%  "}
    @{Row=19; E="  -1.29%  "}
    @{Row=20; E="  -3.08%  "}
    @{Row=21; D="72.46"; E="  +1.54%  "}
    @{Row=22; D="11.14"; E="  -2.52%  "}
    @{Row=23; D="237.47"; E="  -1.25%  "}
    @{Row=24; D="2.73"; E="  -0.22%  "}
    @{Row=25; D="1.93"; E="  +1.20%  "}
    @{Row=26; E="  +0.14%  "}
    @{Row=27; D="24.45"; E="  -1.60%  "}
    @{Row=28; E="  -1.66%  "}
    @{Row=29; D="9.66"; E="  -1.74%  "}
    @{Row=30; D="36.15"; E="  +3.09%  "}
    @{Row=31; D="159.39"; E="  +2.07%  "}
    @{Row=32; E="  -1.76%  "}
    @{Row=33; E="  +0.08%  "}
    @{Row=34; E="  -0.25%  "}
    @{Row=35; D="0.0755"; E="  -1.59%  "}
    @{Row=36; E="  -3.98%  "}
    @{Row=37; B="LidoDAOToken"; C="https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; D="2.90"; E="  +0.47%  "}
    @{Row=38; B="Stellar"; C="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D="0.116"; E="  +0.48%  "}
    @{Row=39; D="1.82"; E="  +1.17%  "}
    @{Row=40; E="  -0.24%  "}
    @{Row=41; D="3.96"; E="  +0.21%  "}
    @{Row=42; E="  -7.17%  "}
    @{Row=43; D="1.990.54"; E="  +0.63%  "}
    @{Row=44; E="  -1.52%  "}
    @{Row=45; D="18.52"; E="  -2.58%  "}
    @{Row=46; E="  -1.81%  "}
    @{Row=47; D="9.68"; E="  +6.43%  "}
    @{Row=48; D="2.696.55"; E="  -0.21%  "}
    @{Row=49; B="Aave"; C="https://coinranking.com/coin/ixgUfzmLR+aave-aave"; D="96.70"; E="  -0.37%  "}
    @{Row=50; B="BitcoinSV"; C="https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"; D="74.56"; E="  +0.47%  "}
    @{Row=51; D="66.73"; E="  -0.03%  "}
)

foreach ($item in $updates) {
    $row = $item.Row

    if ($item.ContainsKey("B")) {
        $ws.Cells.Item($row, 2).Value = $item.B
    }
    if ($item.ContainsKey("C")) {
        $ws.Cells.Item($row, 3).Value = $item.C
    }
    if ($item.ContainsKey("D")) {
        # Force the Price column to be stored as text so values such as
        # "0.548" or "2.90" are not reinterpreted/rounded as numbers.
        $dCell = $ws.Cells.Item($row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $item.D
        $dCell.Style = "Normal"
    }
    if ($item.ContainsKey("E")) {
        $ws.Cells.Item($row, 5).Value = $item.E
    }
}
